$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (67) down onto the two new rows.
[void]$ws.Range("A67:E67").Copy($ws.Range("A68:E68"))
[void]$ws.Range("A67:E67").Copy($ws.Range("A69:E69"))

# The "Runmode" (column D) cell on the new rows uses the plain bordered style
# (same formatting as e.g. B21) instead of the style used by the existing D column
# cells, so copy that cell's format onto D68/D69 specifically.
[void]$ws.Range("B21").Copy($ws.Range("D68"))
[void]$ws.Range("B21").Copy($ws.Range("D69"))

# Fill in the values for the two new test cases (order matches the order the
# strings were first entered, so the shared-string table comes out the same way).
$ws.Range("A68").Value = "Profile67"
$ws.Range("A69").Value = "Profile68"
$ws.Range("B68").Value = "OPQA-2095"
$ws.Range("C68").Value = "Verify that when user clicks outside of the Neon on-boarding welcome modal then the user will be sent right to their desired location in Neon."
$ws.Range("B69").Value = "OPQA-2114"
$ws.Range("C69").Value = "Verify that when user clicks outside of the profile model then user will be sent right to their desired location in Neon."
$ws.Range("D68").Value = "Y"
$ws.Range("D69").Value = "Y"

# Row 68 is tall (wrapped text), like the rows above it; row 69 stays default height.
$ws.Range("A68:E68").RowHeight = 30

# Reflect the resulting scroll/selection state.
$excel.ActiveWindow.ScrollRow = 56
[void]$ws.Range("G67").Select()
